$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells per the diff (Coin/Link/Price/Volume columns).
# Price-column values that look like plain numbers need an explicit
# text format so Excel keeps them as strings (matching the source data,
# which uses "." as a thousands separator, e.g. "43.114.86").

$ws.Cells.Item(2, 4).Value = '43.114.86'
$ws.Cells.Item(2, 5).Value = '  +0.60%  '

$ws.Cells.Item(3, 4).Value = '2.583.86'
$ws.Cells.Item(3, 5).Value = '  +2.28%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.14%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '315.46'
$ws.Cells.Item(5, 5).Value = '  -0.47%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '97.09'
$ws.Cells.Item(6, 5).Value = '  +2.13%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.577'
$ws.Cells.Item(7, 5).Value = '  -0.15%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.996'
$ws.Cells.Item(8, 5).Value = '  -0.45%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.540'
$ws.Cells.Item(9, 5).Value = '  +1.53%  '

$ws.Cells.Item(10, 5).Value = '  -0.73%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0813'
$ws.Cells.Item(11, 5).Value = '  +0.48%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.51'
$ws.Cells.Item(12, 5).Value = '  -0.64%  '

$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.108'
$ws.Cells.Item(13, 5).Value = '  -3.13%  '

$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '2.914.04'
$ws.Cells.Item(14, 5).Value = '  -0.07%  '

$ws.Cells.Item(15, 4).Value = '2.573.05'
$ws.Cells.Item(15, 5).Value = '  +1.78%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '15.25'
$ws.Cells.Item(16, 5).Value = '  -0.01%  '

$ws.Cells.Item(17, 5).Value = '  -0.35%  '

$ws.Cells.Item(18, 4).Value = '43.129.21'
$ws.Cells.Item(18, 5).Value = '  +0.43%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.85'
$ws.Cells.Item(19, 5).Value = '  +3.13%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.54'
$ws.Cells.Item(20, 5).Value = '  -2.88%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0962'
$ws.Cells.Item(21, 5).Value = '  -0.17%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '69.39'
$ws.Cells.Item(22, 5).Value = '  -1.01%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '254.56'
$ws.Cells.Item(23, 5).Value = '  +1.45%  '

$ws.Cells.Item(24, 5).Value = '  +0.78%  '

$ws.Cells.Item(25, 5).Value = '  +3.34%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '27.28'
$ws.Cells.Item(26, 5).Value = '  +1.68%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.999'
$ws.Cells.Item(27, 5).Value = '  -0.09%  '

$ws.Cells.Item(28, 5).Value = '  +1.73%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '40.15'
$ws.Cells.Item(29, 5).Value = '  +1.05%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '10.31'
$ws.Cells.Item(30, 5).Value = '  +0.63%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '5.86'
$ws.Cells.Item(31, 5).Value = '  -2.65%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '155.61'
$ws.Cells.Item(32, 5).Value = '  +0.64%  '

$ws.Cells.Item(33, 5).Value = '  +4.02%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.16'
$ws.Cells.Item(34, 5).Value = '  +2.39%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0807'
$ws.Cells.Item(35, 5).Value = '  +2.25%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.69'
$ws.Cells.Item(36, 5).Value = '  +3.36%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '18.73'
$ws.Cells.Item(37, 5).Value = '  -0.99%  '

$ws.Cells.Item(38, 5).Value = '  -0.10%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.47'
$ws.Cells.Item(39, 5).Value = '  +9.02%  '

$ws.Cells.Item(40, 5).Value = '  -0.41%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '22.49'
$ws.Cells.Item(41, 5).Value = '  -5.38%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.96'
$ws.Cells.Item(42, 5).Value = '  +5.09%  '

$ws.Cells.Item(43, 5).Value = '  +0.17%  '

$ws.Cells.Item(44, 5).Value = '  -0.12%  '

$ws.Cells.Item(45, 5).Value = '  -0.67%  '

$ws.Cells.Item(46, 4).Value = '2.009.87'
$ws.Cells.Item(46, 5).Value = '  -0.44%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '8.95'
$ws.Cells.Item(47, 5).Value = '  +2.01%  '

$ws.Cells.Item(48, 4).Value = '2.815.10'
$ws.Cells.Item(48, 5).Value = '  +1.48%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '82.82'
$ws.Cells.Item(49, 5).Value = '  -3.53%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '75.21'
$ws.Cells.Item(50, 5).Value = '  +2.45%  '

$ws.Cells.Item(51, 5).Value = '  +2.41%  '
